$wb = $excel.ActiveWorkbook

# --- Settings sheet: BaseFolder value gets a clearer placeholder ---
$settings = $wb.Worksheets.Item("Settings")
$settings.Range("B2").Value = "<add root for the project's basefolder>"

# --- Email sheet: ReportRecipient value gets a clearer placeholder,
#     and a new EmailBodyTemplate row is added to the table ---
$email = $wb.Worksheets.Item("Email")
$email.Range("B2").Value = "<add email address to which the report is sent>"

$lo = $email.ListObjects.Item(1)
$newRow = $lo.ListRows.Add()
$email.Range("A3").Value = "EmailBodyTemplate"
$email.Range("B3").Value = "Please find attached travel expense report from today's inputs."
$email.Range("C3").Value = "Email message"

# --- Restore cursor position on Settings, then leave Email as the active tab ---
[void]$settings.Activate()
[void]$settings.Range("C18").Select()

[void]$email.Activate()
[void]$email.Range("E10").Select()
